$d = $word.ActiveDocument

# --- Edit 1: shorten the intro paragraph ---
# Remove the hyperlink wrapper around "PCA" so its text becomes a plain run
# (matches how Word leaves the run behind once the hyperlink field is deleted),
# then remove the now-plain-text sentence referencing PCA/SVM/decision trees.
if ($d.Hyperlinks.Count -gt 0) {
  $d.Hyperlinks(1).Delete()
}

$rng = $d.Content
$found = $rng.Find.Execute(
  "Previously, we managed to implement PCA and next time we will deal with SVM and decision trees.",
  $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
  $rng.Delete()
}

# --- Edit 2: drop the "And this gives us this cool plot" section through the
# "Possible extensions" bullet list at the end of the document, leaving the
# final empty paragraph intact. ---
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $ptext = $d.Paragraphs.Item($i).Range.Text.Trim()
  if ($startPara -eq $null -and $ptext -eq "And this gives us this cool plot:") {
    $startPara = $i
  }
  if ($ptext -like "What is the influence of the number of neighbors *") {
    $endPara = $i
  }
}

if ($startPara -ne $null -and $endPara -ne $null) {
  $p1 = $d.Paragraphs.Item($startPara)
  $p2 = $d.Paragraphs.Item($endPara)
  $delRng = $d.Range($p1.Range.Start, $p2.Range.End)
  $delRng.Delete()
}
